$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.294
$ws.Range("A3").Value = -21.945
$ws.Range("A14").Value = -21.909
$ws.Range("A16").Value = -21.985
$ws.Range("C18").Value = -11.39
$ws.Range("A21").Value = -20.067
$ws.Range("A23").Value = -20.198
$ws.Range("C24").Value = -12.458
$ws.Range("A25").Value = -21.814
$ws.Range("C25").Value = -13.241
$ws.Range("A26").Value = -21.277
$ws.Range("C27").Value = -13.055
$ws.Range("A29").Value = -21.219
$ws.Range("C30").Value = -13.133
$ws.Range("C31").Value = -13.198
$ws.Range("C39").Value = -12.702
$ws.Range("A40").Value = -20
$ws.Range("C42").Value = -12.646
$ws.Range("C48").Value = -11.1
$ws.Range("C51").Value = -11.152
$ws.Range("C52").Value = -11.601
$ws.Range("A53").Value = -21.91
$ws.Range("C55").Value = -13.513
$ws.Range("C56").Value = -13.222
$ws.Range("A57").Value = -22.095
$ws.Range("C57").Value = -13.58
$ws.Range("A59").Value = -22.323
$ws.Range("C60").Value = -12.8
$ws.Range("A65").Value = -21.421
$ws.Range("A69").Value = -21.805
$ws.Range("C73").Value = -12.601
$ws.Range("C74").Value = -12.462
$ws.Range("A79").Value = -21.141
$ws.Range("A83").Value = -21.938
$ws.Range("C89").Value = -10.857
$ws.Range("C90").Value = -12.957
$ws.Range("A91").Value = -21.509
$ws.Range("C92").Value = -11.021
$ws.Range("A93").Value = -21.457
$ws.Range("A100").Value = -21.941
